$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"   = 0.288
    "C5"   = 0.706
    "C6"   = 0.996
    "C7"   = 1.59
    "C8"   = 1.92
    "C10"  = 1.103
    "C11"  = 0.523
    "C12"  = 0.06900000000000001
    "C28"  = 0.098
    "C29"  = 0.238
    "C30"  = 0.369
    "C31"  = 0.376
    "C32"  = 0.289
    "C33"  = 0.201
    "C36"  = 0.018
    "C51"  = 0
    "C52"  = 0.156
    "C53"  = 0.52
    "C54"  = 0.944
    "C55"  = 1.274
    "C56"  = 1.354
    "C57"  = 1.188
    "C58"  = 0.797
    "C78"  = 0.75
    "C79"  = 0.913
    "C80"  = 1.042
    "C81"  = 1.059
    "C100" = 0.138
    "C101" = 0.489
    "C102" = 0.853
    "C103" = 1.252
    "C104" = 1.482
    "C105" = 1.591
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
